$d = $word.ActiveDocument

# 1) Update the header text: "Author / The Title Is Baz / " -> "Author / The Title Is 'Baz' / "
$d.Content.Find.Execute("Author / The Title Is Baz / ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Author / The Title Is ‘Baz’ / ", 2)

# 2) Update the Title paragraph: "Baz" -> "'Baz'" (with curly quotes), as separate runs
$titlePara = $d.Paragraphs.Item(1)
$r = $titlePara.Range
$r.Find.Execute("Baz", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bazStart = $r.Start
$r.Delete()

$p1 = $d.Range($bazStart, $bazStart)
$p1.InsertAfter("‘")

$p2 = $d.Range($bazStart + 1, $bazStart + 1)
$p2.InsertAfter("Baz")

$p3 = $d.Range($bazStart + 4, $bazStart + 4)
$p3.InsertAfter("’")
